# Multi-school management + board tracking: insert 3 new columns
# (School Name, Board, Academic Year) right after "Student ID" (col A),
# shifting all the existing columns (old B..X) to the right by 3 (new E..AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at B:D - this shifts existing data from
# B..X to E..AA, carries header styling/fill along with it, and widens
# the sheet dimension automatically.
$ws.Columns("B:D").Insert()

# Restore / set the explicit widths for the three new columns. Excel
# stores ColumnWidth (character units) internally as
# (OOXML width) = ColumnWidth + 0.8333333... for the default font, so we
# back that padding out to land on the exact target widths (32, 7, 15).
$ws.Columns("B").ColumnWidth = 31.16666666666667
$ws.Columns("C").ColumnWidth = 6.166666666666667
$ws.Columns("D").ColumnWidth = 14.166666666666666

# Header row labels for the new columns.
$ws.Cells.Item(1, 2).Value = "School Name"
$ws.Cells.Item(1, 3).Value = "Board"
$ws.Cells.Item(1, 4).Value = "Academic Year"

# Per-student School Name / Board / Academic Year values, row 2..16.
$rows = @(
  @('Ryan International High School','CBSE','2024-2025'),
  @('Ryan International High School','CBSE','2024-2025'),
  @('Cambridge High School','ICSE','2024-2025'),
  @('Cambridge High School','ICSE','2024-2025'),
  @('St Lawrence High School','CBSE','2024-2025'),
  @('St Lawrence High School','CBSE','2024-2025'),
  @('DAV Public School','CBSE','2024-2025'),
  @('Cathedral School','ICSE','2024-2025'),
  @('Delhi Public School','CBSE','2024-2025'),
  @('Delhi Public School','CBSE','2024-2025'),
  @('Modern High School','ICSE','2024-2025'),
  @('Modern High School','ICSE','2024-2025'),
  @('Ryan International High School','CBSE','2024-2025'),
  @('Cambridge High School','ICSE','2024-2025'),
  @('St Lawrence High School','CBSE','2024-2025')
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
